$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.232.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.663.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("E4").Value = "  +0.45%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5222"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.006"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2647"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.77%  "
$ws.Range("E9").Value = "  -2.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.80"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07718"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.666.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.430"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.891.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5447"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8158"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.253.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.656"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.032"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.05%  "
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.148"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.414"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06112"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.278"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.567"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.250"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.97%  "
$ws.Range("E34").Value = "  -3.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9648"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.58%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.782"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5675"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01600"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.25%  "
$ws.Range("E40").Value = "  -2.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8549"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("E42").Value = "  +0.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.015.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("E45").Value = "  -0.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈112"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.008"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.002"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.482"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05183"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.50%  "
